$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.472.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6307"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07683"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2935"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07748"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.873.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.037"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("E14").Value = "  +0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001067"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.121.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.201"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.488.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.470"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1384"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.416"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("E29").Value = "  +5.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.470"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.87%  "

$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.136"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.044"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.855"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.27%  "

$ws.Range("E35").Value = "  +1.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7090"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.587"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.784"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01793"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.220.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.553"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9097"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.029.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000120"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.132"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "

$ws.Range("E49").Value = "  +1.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.042"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.33%  "

$ws.Range("E51").Value = "  +0.97%  "
